# Scheduled runner update: refresh market-price derived columns (H-N) across leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1020.6
$ws.Range("I29").Value = 34.333332
$ws.Range("J29").Value = 2500
$ws.Range("K29").Value = 102.999996
$ws.Range("L29").Value = 7500
$ws.Range("M29").Value = 178.000004
$ws.Range("N29").Value = -8062
$ws.Range("H34").Value = 1362.5
$ws.Range("I34").Value = 1362.5
$ws.Range("K34").Value = 1362.5
$ws.Range("M34").Value = -1159.5
$ws.Range("H36").Value = 1362.5
$ws.Range("I36").Value = 1362.5
$ws.Range("K36").Value = 1362.5
$ws.Range("M36").Value = -647.5
$ws.Range("H38").Value = 113.5
$ws.Range("I38").Value = 113.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 340.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 31.5
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 1521.1428
$ws.Range("I40").Value = 1629.6
$ws.Range("K40").Value = 1629.6
$ws.Range("M40").Value = -1454.6
$ws.Range("H53").Value = 121
$ws.Range("I53").Value = 109.5
$ws.Range("K53").Value = 109.5
$ws.Range("M53").Value = 527.5
$ws.Range("H58").Value = 2123.7144
$ws.Range("J58").Value = 3339
$ws.Range("L58").Value = 10017
$ws.Range("N58").Value = -10317
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H100").Value = 2459
$ws.Range("I100").Value = 2461.25
$ws.Range("K100").Value = 2461.25
$ws.Range("M100").Value = -1920.25
$ws.Range("H141").Value = 1619.3334
$ws.Range("I141").Value = 1619.3334
$ws.Range("K141").Value = 4858.0002
$ws.Range("M141").Value = 321.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1395.6825
$ws.Range("I32").Value = 1219.305
$ws.Range("K32").Value = 1219.305
$ws.Range("M32").Value = -932.3050000000001
$ws.Range("H45").Value = 2498
$ws.Range("I45").Value = 2496
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 2496
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -2119
$ws.Range("N45").Value = -3254
$ws.Range("H102").Value = 919.7857
$ws.Range("I102").Value = 729
$ws.Range("J102").Value = 3400
$ws.Range("K102").Value = 729
$ws.Range("L102").Value = 3400
$ws.Range("M102").Value = 893
$ws.Range("N102").Value = -6644

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4697.857
$ws.Range("I20").Value = 4314.5
$ws.Range("J20").Value = 6998
$ws.Range("K20").Value = 4314.5
$ws.Range("L20").Value = 6998
$ws.Range("M20").Value = -4067.5
$ws.Range("N20").Value = -7492
$ws.Range("H105").Value = 3308.7
$ws.Range("I105").Value = 3336
$ws.Range("K105").Value = 3336
$ws.Range("M105").Value = -1589
$ws.Range("H107").Value = 2587.5
$ws.Range("I107").Value = 5000
$ws.Range("K107").Value = 5000
$ws.Range("M107").Value = -3080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1800.75
$ws.Range("I6").Value = 825
$ws.Range("K6").Value = 825
$ws.Range("M6").Value = -712
$ws.Range("H17").Value = 12987.333
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4826
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H23").Value = 55505
$ws.Range("J23").Value = 55505
$ws.Range("L23").Value = 55505
$ws.Range("N23").Value = -55985
$ws.Range("H25").Value = 3585.5
$ws.Range("I25").Value = 447.33334
$ws.Range("J25").Value = 13000
$ws.Range("K25").Value = 447.33334
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = -273.33334
$ws.Range("N25").Value = -13348
$ws.Range("H27").Value = 55505
$ws.Range("J27").Value = 55505
$ws.Range("L27").Value = 55505
$ws.Range("N27").Value = -55889
$ws.Range("H41").Value = 18499.125
$ws.Range("J41").Value = 23332.334
$ws.Range("L41").Value = 23332.334
$ws.Range("N41").Value = -24188.334
$ws.Range("H59").Value = 34937
$ws.Range("J59").Value = 34937
$ws.Range("L59").Value = 34937
$ws.Range("N59").Value = -37227
$ws.Range("H60").Value = 24855.857
$ws.Range("J60").Value = 24998.5
$ws.Range("L60").Value = 24998.5
$ws.Range("N60").Value = -26020.5
$ws.Range("H68").Value = 39997.35
$ws.Range("J68").Value = 39997.35
$ws.Range("L68").Value = 39997.35
$ws.Range("N68").Value = -41495.35
$ws.Range("H71").Value = 39997.35
$ws.Range("J71").Value = 39997.35
$ws.Range("L71").Value = 119992.05
$ws.Range("N71").Value = -127480.05
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2540
$ws.Range("I5").Value = 3458
$ws.Range("K5").Value = 10374
$ws.Range("M5").Value = -10262
$ws.Range("H39").Value = 5049.25
$ws.Range("I39").Value = 1649.5
$ws.Range("J39").Value = 8449
$ws.Range("K39").Value = 4948.5
$ws.Range("L39").Value = 25347
$ws.Range("M39").Value = -4654.5
$ws.Range("N39").Value = -25935
$ws.Range("H55").Value = 1696.4286
$ws.Range("I55").Value = 598.3333
$ws.Range("K55").Value = 1794.9999
$ws.Range("M55").Value = -1617.9999
$ws.Range("H135").Value = 2540
$ws.Range("I135").Value = 3458
$ws.Range("K135").Value = 31122
$ws.Range("M135").Value = -28587

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 16706.25
$ws.Range("J46").Value = 18021.428
$ws.Range("L46").Value = 18021.428
$ws.Range("N46").Value = -18333.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 342301.34
$ws.Range("J19").Value = 13452
$ws.Range("L19").Value = 13452
$ws.Range("N19").Value = -13792
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H46").Value = 2943.8
$ws.Range("J46").Value = 3261.9048
$ws.Range("L46").Value = 3261.9048
$ws.Range("N46").Value = -3637.9048
$ws.Range("H68").Value = 36235
$ws.Range("I68").Value = 1956
$ws.Range("K68").Value = 1956
$ws.Range("M68").Value = -1207
$ws.Range("H71").Value = 36235
$ws.Range("I71").Value = 1956
$ws.Range("K71").Value = 9780
$ws.Range("M71").Value = -6036
$ws.Range("H100").Value = 1600
$ws.Range("I100").Value = 1600
$ws.Range("K100").Value = 1600
$ws.Range("M100").Value = -1059
$ws.Range("H122").Value = 3504
$ws.Range("I122").Value = 3504
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10512
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8062
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 18028
$ws.Range("I54").Value = 70
$ws.Range("K54").Value = 70
$ws.Range("M54").Value = 450
$ws.Range("H96").Value = 2911.75
$ws.Range("I96").Value = 2401
$ws.Range("J96").Value = 4444
$ws.Range("K96").Value = 2401
$ws.Range("L96").Value = 4444
$ws.Range("M96").Value = -1028
$ws.Range("N96").Value = -7190
